$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (Wins, Losses, Ties) after the existing
# "Unnamed: 28" column (AC), matching the style of the other header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record (W/L/T) for every data row.
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 79  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 83  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
